$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 285, pushing ADL..PPT down by one row (to 286..297)
$ws.Rows.Item(285).Insert()

# Copy the style of the colo-code column (A) from the row below (now ADL at 286)
# onto the newly inserted row's A cell, so it keeps the bold/bordered/centered style.
$ws.Cells.Item(286, 1).Copy()
$ws.Cells.Item(285, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 285 with the Kingston, Jamaica colo entry
$ws.Cells.Item(285, 1).Value = "KIN"
$ws.Cells.Item(285, 2).Value = "Kingston"
$ws.Cells.Item(285, 3).Value = 17.9951
$ws.Cells.Item(285, 4).Value = -76.7846
$ws.Cells.Item(285, 5).Value = "JM"
$ws.Cells.Item(285, 6).Value = "North America"
$ws.Cells.Item(285, 7).Value = "Kingston"
